# Apply updated dSF (column F) values for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    7  = 0
    10 = 1
    14 = 0
    23 = -5
    25 = -4
    26 = 0
    28 = -6
    36 = -6
    37 = -3
    38 = -1
    46 = -1
    47 = -5
    48 = 7
    49 = -3
    59 = -4
    62 = 0
    67 = -8
    69 = 2
    72 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
